$wb = $excel.ActiveWorkbook

# 1. Rename sheet "version" to "other_questions"
$verSheet = $wb.Worksheets.Item("version")
$verSheet.Name = "other_questions"

# 2. Update the Bernabeu question text on the questions sheet
$qWs = $wb.Worksheets.Item("questions")
$qWs.Range("C25").Value = '¿Cuántos estadios como el Bernabéu podrías cubrir con todo el papel y cartón que usamos desde los hogares de los españoles?'

# 3. Update other_questions sheet headers and data
$ws = $wb.Worksheets.Item("other_questions")
$ws.Range("A1").Value = "id_play_question"
$ws.Range("B1").Value = "question"
$ws.Columns.Item(1).ColumnWidth = 16

$questions = @(
    '¿Cuántos folios A4 se necesitan para cubrir por completo la superficie de un campo de fútbol?'
    '¿Cuántas gotas de agua hay en una botella estándar de 1,5 litros?'
    '¿Cuántas veces parpadea una persona de media a lo largo de toda su vida?'
    '¿Cuántas hormigas viven actualmente en toda la península ibérica?'
    '¿Cuántos pelos tiene un gato doméstico de tamaño medio?'
    '¿Cuántos litros de combustible consume un Airbus A380 en un vuelo de Madrid a Nueva York?'
    '¿Cuántas pizzas se reparten en Madrid un sábado por la noche entre las 20:00 y las 23:00?'
    '¿Cuántas bombillas hay instaladas en el total de los hogares de España?'
    '¿Cuántas palabras pronuncia una persona media a lo largo de un día normal?'
    '¿Cuántos granos de arena hay en la playa de la Malvarrosa en Valencia?'
    'Uno de cada _________ españoles ha nacido un 29 de febrero. '
    '¿Cuántas calorías quema un corredor de maratón durante los 42 kilómetros de carrera?'
    '¿Cuántos mensajes de WhatsApp se envían en todo el mundo durante un solo segundo?'
    '¿Cuántas manzanas de tamaño medio se necesitan para llenar el maletero de un coche utilitario?'
    '¿Cuántos litros de leche produce una vaca a lo largo de toda su vida productiva?'
    '¿Cuántas hormigas comunes serían necesarias para igualar el peso de un elefante africano?'
    '¿Cuántos neumáticos de coche se desechan en España cada año?'
    '¿Cuántos minutos de video se suben a YouTube cada hora a nivel global?'
    '¿Cuántas monedas de 5 céntimos harían falta para pagar el fichaje de un jugador de fútbol de 100 millones de euros?'
    '¿Cuántos kilómetros totales recorre una abeja para producir un solo kilo de miel?'
    '¿Cuántos árboles hay plantados en el Parque del Retiro de Madrid?'
    '¿Cuántas patatas fritas de bolsa se consumen en los cines de España durante el fin de semana de un gran estreno?'
    '¿Cuántas latas de refresco se reciclan en Europa en un día cualquiera?'
    '¿Cuántas células tiene el cuerpo de un bebé recién nacido?'
    'Uno de cada _________ habitantes del mundo vive actualmente en China o en la India. '
    '¿Cuántas gotas de lluvia caen sobre una ciudad del tamaño de Sevilla durante una tormenta de 10 minutos?'
    '¿Cuántas veces cabe el volumen de la Luna dentro del volumen de la Tierra?'
    '¿Cuántos cigarrillos se fuman en España en un solo día de diario?'
    '¿Cuántos metros de hilo dental usa una persona media si se limpia los dientes cada noche de su vida?'
    '¿Cuántos globos de helio estándar se necesitarían para levantar del suelo a una persona de 70 kg?'
    '¿Cuántas horas de sueño acumuladas tiene una persona al cumplir los 80 años?'
    '¿Cuántas piezas de LEGO existen actualmente en el mundo (fabricadas desde su inicio)?'
    '¿Cuántas veces late el corazón de un ratón de campo en un solo minuto?'
    '¿Cuántos litros de agua se evaporan de una piscina olímpica en un día caluroso de agosto?'
    '¿Cuántas fotos digitales hay guardadas en la nube por cada habitante del planeta?'
    '¿Cuántas baldosas hay en las aceras de una calle principal de 1 kilómetro de largo?'
    '¿Cuántas abejas viven en una colmena típica durante la primavera?'
    '¿Cuántas tazas de té se beben en el Reino Unido durante una tarde de domingo?'
    '¿Cuántos pasos hay que dar de media para subir al Everest partiendo desde el campamento base?'
    '¿Cuántos tornillos y remaches tiene un avión Boeing 747?'
    '¿Cuántas búsquedas en Google se realizan desde España durante los 90 minutos de una final de Champions?'
    '¿Cuántas tarjetas de crédito caben (en volumen) dentro de una cabina telefónica clásica?'
    '¿Cuántos átomos de hidrógeno hay en un vaso de agua de 200 ml?'
    '¿Cuántos kilos de basura genera un festival de música multitudinario durante tres días?'
    '¿Cuántas veces se podría escribir "El Quijote" completo con la tinta de un bolígrafo BIC convencional?'
    '¿Cuántas moscas hay en una granja de vacas de tamaño medio en verano?'
    'Uno de cada _________ coches que circulan por Madrid es de color blanco. '
    '¿Cuántos kilómetros de vasos de plástico se formarían si pusiéramos en fila todos los usados en un gran concierto?'
    '¿Cuántas respiraciones da un perro de raza media a lo largo de sus 12 años de vida?'
    '¿Cuántos paneles solares se necesitarían para cubrir toda la demanda eléctrica de España?'
    '¿Cuántos litros de pintura se necesitaron para pintar el exterior de la Casa Blanca?'
    '¿Cuántas personas están volando en un avión en este preciso instante en todo el planeta?'
    '¿Cuántas palomitas de maíz caben en el habitáculo de un autobús escolar?'
    '¿Cuántos granos de azúcar hay en un sobre individual de 8 gramos?'
    '¿Cuántos metros crece el pelo de toda la población de una ciudad de 100.000 habitantes en un mes?'
    '¿Cuántas llaves de casa hay perdidas u olvidadas en el fondo del mar?'
    '¿Cuántas veces se ha reproducido el video más visto de YouTube por cada persona viva en la Tierra?'
    '¿Cuántas veces daría la vuelta al mundo el cableado eléctrico de una ciudad como Barcelona?'
    '¿Cuántos libros físicos hay en todas las bibliotecas públicas de España juntas?'
    '¿Cuántas hamburguesas vende McDonald''s en todo el mundo en un solo minuto?'
)

$lastRow = 1 + $questions.Length

# Copy the existing formatted row 2 down across the new rows so every new cell
# keeps the same style as the original data row before we overwrite the values.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $questions.Length; $i++) {
    $row = $i + 2
    $id = 101 + $i
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $questions[$i]
}

Write-Host "done"